$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 14 de Junio de 2020 a las 17:40"

# Update country data rows (country name shuffles + refreshed COVID-19 stats)
# Row 4: Estados Unidos -> Estados Unidos
$ws.Cells.Item(4, 2).Value = 2146884
$ws.Cells.Item(4, 3).Value = 4660
$ws.Cells.Item(4, 5).Value = 1175219
$ws.Cells.Item(4, 7).Value = 32
$ws.Cells.Item(4, 8).Value = 117559

# Row 7: India -> India
$ws.Cells.Item(7, 4).Value = 164803
$ws.Cells.Item(7, 5).Value = 150432

# Row 8: Reino Unido -> Reino Unido
$ws.Cells.Item(8, 2).Value = 295889
$ws.Cells.Item(8, 3).Value = 1514
$ws.Cells.Item(8, 7).Value = 36
$ws.Cells.Item(8, 8).Value = 41698

# Row 12: Iran -> Alemania
$ws.Cells.Item(12, 1).Value = "Alemania"
$ws.Cells.Item(12, 2).Value = 187489
$ws.Cells.Item(12, 3).Value = 66
$ws.Cells.Item(12, 4).Value = 172200
$ws.Cells.Item(12, 5).Value = 6422
$ws.Cells.Item(12, 7).Value = 0
$ws.Cells.Item(12, 8).Value = 8867

# Row 13: Alemania -> Iran
$ws.Cells.Item(13, 1).Value = "Iran"
$ws.Cells.Item(13, 2).Value = 187427
$ws.Cells.Item(13, 3).Value = 2472
$ws.Cells.Item(13, 4).Value = 148674
$ws.Cells.Item(13, 5).Value = 29916
$ws.Cells.Item(13, 7).Value = 107
$ws.Cells.Item(13, 8).Value = 8837

# Row 20: Canada -> Canada
$ws.Cells.Item(20, 2).Value = 98607
$ws.Cells.Item(20, 3).Value = 197
$ws.Cells.Item(20, 4).Value = 59777
$ws.Cells.Item(20, 5).Value = 30711
$ws.Cells.Item(20, 7).Value = 12
$ws.Cells.Item(20, 8).Value = 8119

# Row 33: Singapur -> Singapur
$ws.Cells.Item(33, 4).Value = 29589
$ws.Cells.Item(33, 5).Value = 10989

# Row 45: Republica Dominicana -> Republica Dominicana
$ws.Cells.Item(45, 2).Value = 22962
$ws.Cells.Item(45, 3).Value = 390
$ws.Cells.Item(45, 4).Value = 13320
$ws.Cells.Item(45, 5).Value = 9050
$ws.Cells.Item(45, 7).Value = 15
$ws.Cells.Item(45, 8).Value = 592

# Row 47: Panama -> Irak
$ws.Cells.Item(47, 1).Value = "Irak"
$ws.Cells.Item(47, 2).Value = 20209
$ws.Cells.Item(47, 3).Value = 1259
$ws.Cells.Item(47, 4).Value = 8121
$ws.Cells.Item(47, 5).Value = 11481
$ws.Cells.Item(47, 7).Value = 58
$ws.Cells.Item(47, 8).Value = 607

# Row 48: Israel -> Panama
$ws.Cells.Item(48, 1).Value = "Panama"
$ws.Cells.Item(48, 2).Value = 20059
$ws.Cells.Item(48, 3).Value = 0
$ws.Cells.Item(48, 4).Value = 13759
$ws.Cells.Item(48, 5).Value = 5871
$ws.Cells.Item(48, 8).Value = 429

# Row 49: Irak -> Israel
$ws.Cells.Item(49, 1).Value = "Israel"
$ws.Cells.Item(49, 2).Value = 19008
$ws.Cells.Item(49, 3).Value = 36
$ws.Cells.Item(49, 4).Value = 15360
$ws.Cells.Item(49, 5).Value = 3348
$ws.Cells.Item(49, 8).Value = 300

# Row 56: Kazajistan -> Kazajistan
$ws.Cells.Item(56, 4).Value = 9174
$ws.Cells.Item(56, 5).Value = 5249

# Row 63: Chequia -> Chequia
$ws.Cells.Item(63, 2).Value = 9999
$ws.Cells.Item(63, 3).Value = 8
$ws.Cells.Item(63, 5).Value = 2451
$ws.Cells.Item(63, 7).Value = 1
$ws.Cells.Item(63, 8).Value = 329

# Row 64: Azerbaiyan -> Azerbaiyan
$ws.Cells.Item(64, 2).Value = 9957
$ws.Cells.Item(64, 3).Value = 387
$ws.Cells.Item(64, 4).Value = 5583
$ws.Cells.Item(64, 5).Value = 4255
$ws.Cells.Item(64, 7).Value = 4
$ws.Cells.Item(64, 8).Value = 119

# Row 77: Tayikistan -> Tayikistan
$ws.Cells.Item(77, 2).Value = 5035
$ws.Cells.Item(77, 3).Value = 64
$ws.Cells.Item(77, 4).Value = 3409
$ws.Cells.Item(77, 5).Value = 1576

# Row 81: Republica de Yibuti -> Republica de Yibuti
$ws.Cells.Item(81, 2).Value = 4465
$ws.Cells.Item(81, 3).Value = 16
$ws.Cells.Item(81, 4).Value = 2950
$ws.Cells.Item(81, 5).Value = 1472
$ws.Cells.Item(81, 7).Value = 2
$ws.Cells.Item(81, 8).Value = 43

# Row 83: Hungria -> Luxemburgo
$ws.Cells.Item(83, 1).Value = "Luxemburgo"
$ws.Cells.Item(83, 2).Value = 4070
$ws.Cells.Item(83, 3).Value = 7
$ws.Cells.Item(83, 4).Value = 3929
$ws.Cells.Item(83, 5).Value = 31
$ws.Cells.Item(83, 7).Value = 0
$ws.Cells.Item(83, 8).Value = 110

# Row 84: Luxemburgo -> Hungria
$ws.Cells.Item(84, 1).Value = "Hungria"
$ws.Cells.Item(84, 2).Value = 4069
$ws.Cells.Item(84, 3).Value = 5
$ws.Cells.Item(84, 4).Value = 2482
$ws.Cells.Item(84, 5).Value = 1025
$ws.Cells.Item(84, 7).Value = 3
$ws.Cells.Item(84, 8).Value = 562

# Row 85: Republica de Macedonia -> Republica de Macedonia
$ws.Cells.Item(85, 2).Value = 4057
$ws.Cells.Item(85, 3).Value = 162
$ws.Cells.Item(85, 4).Value = 1710
$ws.Cells.Item(85, 5).Value = 2159
$ws.Cells.Item(85, 7).Value = 9
$ws.Cells.Item(85, 8).Value = 188

# Row 92: Grecia -> Grecia
$ws.Cells.Item(92, 2).Value = 3121
$ws.Cells.Item(92, 3).Value = 9
$ws.Cells.Item(92, 5).Value = 1564

# Row 121: Sierra Leona -> Guayana Francesa
$ws.Cells.Item(121, 1).Value = "Guayana Francesa"
$ws.Cells.Item(121, 2).Value = 1255
$ws.Cells.Item(121, 3).Value = 94
$ws.Cells.Item(121, 4).Value = 534
$ws.Cells.Item(121, 5).Value = 718
$ws.Cells.Item(121, 7).Value = 1
$ws.Cells.Item(121, 8).Value = 3

# Row 122: Guayana Francesa -> Sierra Leona
$ws.Cells.Item(122, 1).Value = "Sierra Leona"
$ws.Cells.Item(122, 2).Value = 1169
$ws.Cells.Item(122, 3).Value = 37
$ws.Cells.Item(122, 4).Value = 680
$ws.Cells.Item(122, 5).Value = 438
$ws.Cells.Item(122, 8).Value = 51

# Row 128: Jordania -> Jordania
$ws.Cells.Item(128, 2).Value = 961
$ws.Cells.Item(128, 3).Value = 8
$ws.Cells.Item(128, 4).Value = 682
$ws.Cells.Item(128, 5).Value = 270

# Row 134: Congo -> Cabo Verde
$ws.Cells.Item(134, 1).Value = "Cabo Verde"
$ws.Cells.Item(134, 2).Value = 750
$ws.Cells.Item(134, 3).Value = 24
$ws.Cells.Item(134, 4).Value = 301
$ws.Cells.Item(134, 5).Value = 443
$ws.Cells.Item(134, 8).Value = 6

# Row 135: Cabo Verde -> Congo
$ws.Cells.Item(135, 1).Value = "Congo"
$ws.Cells.Item(135, 2).Value = 728
$ws.Cells.Item(135, 4).Value = 221
$ws.Cells.Item(135, 5).Value = 483
$ws.Cells.Item(135, 8).Value = 24

# Row 148: Estado de Palestina -> Reunion
$ws.Cells.Item(148, 1).Value = "Reunion"
$ws.Cells.Item(148, 2).Value = 495
$ws.Cells.Item(148, 3).Value = 6
$ws.Cells.Item(148, 4).Value = 460
$ws.Cells.Item(148, 5).Value = 34
$ws.Cells.Item(148, 8).Value = 1

# Row 149: Reunion -> Estado de Palestina
$ws.Cells.Item(149, 1).Value = "Estado de Palestina"
$ws.Cells.Item(149, 4).Value = 415
$ws.Cells.Item(149, 5).Value = 71
$ws.Cells.Item(149, 8).Value = 3

# Row 206: Groenlandia -> Islas Malvinas
$ws.Cells.Item(206, 1).Value = "Islas Malvinas"

# Row 207: Islas Malvinas -> Groenlandia
$ws.Cells.Item(207, 1).Value = "Groenlandia"

# Row 208: Islas Turcas y Caicos -> Santa Sede
$ws.Cells.Item(208, 1).Value = "Santa Sede"
$ws.Cells.Item(208, 4).Value = 12
$ws.Cells.Item(208, 8).Value = 0

# Row 209: Santa Sede -> Islas Turcas y Caicos
$ws.Cells.Item(209, 1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(209, 4).Value = 11
$ws.Cells.Item(209, 8).Value = 1

# Row 213: Papua Nueva Guinea -> Islas Virgenes Britanicas
$ws.Cells.Item(213, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(213, 4).Value = 7
$ws.Cells.Item(213, 8).Value = 1

# Row 214: Islas Virgenes Britanicas -> Papua Nueva Guinea
$ws.Cells.Item(214, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(214, 4).Value = 8
$ws.Cells.Item(214, 8).Value = 0
